$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / safe string assignments (Coin name, Link URL, Volume % text, and
#     Price values whose format is never auto-parsed as a pure number) ---
$ws.Range('D2').Value = '95.030.68'
$ws.Range('E2').Value = '  -1.87%  '
$ws.Range('D3').Value = '3.581.41'
$ws.Range('E3').Value = '  -2.83%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('E5').Value = '  +17.92%  '
$ws.Range('E6').Value = '  -4.91%  '
$ws.Range('E7').Value = '  -3.10%  '
$ws.Range('E8').Value = '  -3.61%  '
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '3.581.23'
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('E12').Value = '  +4.27%  '
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('E14').Value = '  -4.72%  '
$ws.Range('E15').Value = '  -4.84%  '
$ws.Range('D16').Value = '4.275.30'
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('D17').Value = '94.844.09'
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('E18').Value = '  -3.40%  '
$ws.Range('E19').Value = '  +6.64%  '
$ws.Range('D20').Value = '3.585.39'
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('E23').Value = '  -4.24%  '
$ws.Range('E24').Value = '  -5.62%  '
$ws.Range('E25').Value = '  +15.89%  '
$ws.Range('E26').Value = '  +15.19%  '
$ws.Range('E27').Value = '  -4.33%  '
$ws.Range('E28').Value = '  -2.86%  '
$ws.Range('D29').Value = '3.790.80'
$ws.Range('E29').Value = '  -2.42%  '
$ws.Range('E30').Value = '  -6.31%  '
$ws.Range('E31').Value = '  +2.21%  '
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -5.82%  '
$ws.Range('E36').Value = '  -5.49%  '
$ws.Range('E37').Value = '  -2.30%  '
$ws.Range('E38').Value = '  -1.85%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E40').Value = '  -8.99%  '
$ws.Range('E41').Value = '  -6.42%  '
$ws.Range('E42').Value = '  -1.59%  '
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E46').Value = '  +2.64%  '
$ws.Range('B47').Value = 'ImmutableX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E47').Value = '  -7.32%  '
$ws.Range('E48').Value = '  -4.02%  '
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('B50').Value = 'MantraDAO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('E50').Value = '  +1.85%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E51').Value = '  -0.70%  '

# --- Price values that look like plain numbers: Excel would silently coerce a normal
#     Value assignment into a Number (dropping the original text formatting, e.g.
#     "1.00" -> 1). Enter them as a text-producing formula, then Copy/PasteSpecial
#     (values-only) to collapse the formula down to a literal string cell, matching the
#     original inline-string storage exactly without touching cell style/number format. ---
$ws.Range('D4').Formula = '="0.999"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('D6').Formula = '="224.69"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('D7').Formula = '="634.46"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('D9').Formula = '="1.07"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('D12').Formula = '="45.94"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('D14').Formula = '="0.0000286"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('D15').Formula = '="6.41"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('D18').Formula = '="8.72"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('D19').Formula = '="19.90"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('D21').Formula = '="12.81"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('D22').Formula = '="0.508"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('D23').Formula = '="499.29"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('D26').Formula = '="116.39"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('D28').Formula = '="6.70"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('D30').Formula = '="12.47"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('D32').Formula = '="2.87"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('D33').Formula = '="0.999"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('D34').Formula = '="1.00"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('D35').Formula = '="0.177"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('D36').Formula = '="1.75"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('D37').Formula = '="31.40"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('D38').Formula = '="0.579"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('D40').Formula = '="586.01"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('D41').Formula = '="8.23"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('D42').Formula = '="6.72"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('D43').Formula = '="40.34"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('D45').Formula = '="0.465"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('D46').Formula = '="0.0468"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('D47').Formula = '="1.89"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('D48').Formula = '="0.912"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('D49').Formula = '="23.41"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="3.59"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('D51').Formula = '="8.46"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
